$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 468.65216
$ws.Range("I33").Value = 369.41666
$ws.Range("K33").Value = 369.41666
$ws.Range("M33").Value = -140.41666
$ws.Range("H86").Value = 5077.7144
$ws.Range("I86").Value = 4845.357
$ws.Range("J86").Value = 5542.4287
$ws.Range("K86").Value = 4845.357
$ws.Range("L86").Value = 5542.4287
$ws.Range("M86").Value = -3722.357
$ws.Range("N86").Value = -7788.4287
$ws.Range("H89").Value = 5077.7144
$ws.Range("I89").Value = 4845.357
$ws.Range("J89").Value = 5542.4287
$ws.Range("K89").Value = 24226.785
$ws.Range("L89").Value = 27712.1435
$ws.Range("M89").Value = -18610.785
$ws.Range("N89").Value = -38944.14350000001
$ws.Range("H136").Value = 92500
$ws.Range("J136").Value = 92500
$ws.Range("L136").Value = 92500
$ws.Range("N136").Value = -102700
$ws.Range("H137").Value = 97190.734
$ws.Range("I137").Value = 128904
$ws.Range("J137").Value = 8393.6
$ws.Range("K137").Value = 386712
$ws.Range("L137").Value = 25180.8
$ws.Range("M137").Value = -384162
$ws.Range("N137").Value = -30280.8
$ws.Range("H138").Value = 7263.846
$ws.Range("I138").Value = 6498
$ws.Range("J138").Value = 7288.1587
$ws.Range("K138").Value = 19494
$ws.Range("L138").Value = 21864.4761
$ws.Range("M138").Value = -14354
$ws.Range("N138").Value = -32144.4761

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15874.349
$ws.Range("I32").Value = 12827.459
$ws.Range("J32").Value = 34663.5
$ws.Range("K32").Value = 12827.459
$ws.Range("L32").Value = 34663.5
$ws.Range("M32").Value = -12540.459
$ws.Range("N32").Value = -35237.5
$ws.Range("H61").Value = 5903.75
$ws.Range("I61").Value = 5871.7085
$ws.Range("K61").Value = 5871.7085
$ws.Range("M61").Value = -5659.7085
$ws.Range("H74").Value = 85355
$ws.Range("I74").Value = 5001.7
$ws.Range("J74").Value = 888888
$ws.Range("K74").Value = 5001.7
$ws.Range("L74").Value = 888888
$ws.Range("M74").Value = -4127.7
$ws.Range("N74").Value = -890636
$ws.Range("H77").Value = 85355
$ws.Range("I77").Value = 5001.7
$ws.Range("J77").Value = 888888
$ws.Range("K77").Value = 25008.5
$ws.Range("L77").Value = 4444440
$ws.Range("M77").Value = -20640.5
$ws.Range("N77").Value = -4453176
$ws.Range("H88").Value = 3812.6155
$ws.Range("J88").Value = 2956
$ws.Range("L88").Value = 2956
$ws.Range("N88").Value = -3768
$ws.Range("H91").Value = 3812.6155
$ws.Range("J91").Value = 2956
$ws.Range("L91").Value = 2956
$ws.Range("N91").Value = -5764
$ws.Range("H101").Value = 49997.5
$ws.Range("J101").Value = 49997.5
$ws.Range("L101").Value = 49997.5
$ws.Range("N101").Value = -56487.5
$ws.Range("H102").Value = 3975932
$ws.Range("I102").Value = 3975932
$ws.Range("K102").Value = 3975932
$ws.Range("M102").Value = -3974310
$ws.Range("H132").Value = 34776.742
$ws.Range("I132").Value = 2699.5
$ws.Range("K132").Value = 8098.5
$ws.Range("M132").Value = -5568.5
$ws.Range("H136").Value = 5903.75
$ws.Range("I136").Value = 5871.7085
$ws.Range("K136").Value = 17615.1255
$ws.Range("M136").Value = -15065.1255

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2361782
$ws.Range("I94").Value = 3677881.8
$ws.Range("J94").Value = 6656.3687
$ws.Range("K94").Value = 3677881.8
$ws.Range("L94").Value = 6656.3687
$ws.Range("M94").Value = -3677430.8
$ws.Range("N94").Value = -7558.3687
$ws.Range("H107").Value = 23812772
$ws.Range("I107").Value = 23812772
$ws.Range("K107").Value = 23812772
$ws.Range("M107").Value = -23810852
$ws.Range("H134").Value = 8732.5
$ws.Range("I134").Value = 1843.3334
$ws.Range("K134").Value = 5530.0002
$ws.Range("M134").Value = -2995.0002

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 90910590
$ws.Range("I22").Value = 1955.25
$ws.Range("J22").Value = 333333630
$ws.Range("K22").Value = 1955.25
$ws.Range("L22").Value = 333333630
$ws.Range("M22").Value = -1605.25
$ws.Range("N22").Value = -333334330
$ws.Range("H31").Value = 21855.857
$ws.Range("I31").Value = 2991
$ws.Range("K31").Value = 2991
$ws.Range("M31").Value = -2696
$ws.Range("H34").Value = 21855.857
$ws.Range("I34").Value = 2991
$ws.Range("K34").Value = 2991
$ws.Range("M34").Value = -2789
$ws.Range("H132").Value = 92376.71000000001
$ws.Range("I132").Value = 60971.766
$ws.Range("K132").Value = 182915.298
$ws.Range("M132").Value = -180385.298
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140
$ws.Range("H141").Value = 379165.66
$ws.Range("J141").Value = 445998.8
$ws.Range("L141").Value = 445998.8
$ws.Range("N141").Value = -456358.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2754.4546
$ws.Range("I3").Value = 1900
$ws.Range("J3").Value = 3779.8
$ws.Range("K3").Value = 5700
$ws.Range("L3").Value = 11339.4
$ws.Range("M3").Value = -5588
$ws.Range("N3").Value = -11563.4

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 7266.8
$ws.Range("I107").Value = 8920.75
$ws.Range("J107").Value = 651
$ws.Range("K107").Value = 8920.75
$ws.Range("L107").Value = 651
$ws.Range("M107").Value = -7000.75
$ws.Range("N107").Value = -4491
$ws.Range("H122").Value = 3405.9614
$ws.Range("I122").Value = 3067.652
$ws.Range("K122").Value = 9202.956
$ws.Range("M122").Value = -6752.956
$ws.Range("H126").Value = 5369255.5
$ws.Range("I126").Value = 3249328.2
$ws.Range("J126").Value = 8337153.5
$ws.Range("K126").Value = 9747984.600000001
$ws.Range("L126").Value = 25011460.5
$ws.Range("M126").Value = -9745514.600000001
$ws.Range("N126").Value = -25016400.5
$ws.Range("H132").Value = 4261
$ws.Range("I132").Value = 4294.125
$ws.Range("K132").Value = 12882.375
$ws.Range("M132").Value = -10352.375
$ws.Range("H141").Value = 116666.664
$ws.Range("J141").Value = 116666.664
$ws.Range("L141").Value = 116666.664
$ws.Range("N141").Value = -127026.664

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1478714.4
$ws.Range("I2").Value = 10000000
$ws.Range("J2").Value = 58500.168
$ws.Range("K2").Value = 10000000
$ws.Range("L2").Value = 58500.168
$ws.Range("M2").Value = -9999888
$ws.Range("N2").Value = -58724.168
$ws.Range("H40").Value = 7810
$ws.Range("I40").Value = 5030.5713
$ws.Range("K40").Value = 5030.5713
$ws.Range("M40").Value = -4894.5713
$ws.Range("H93").Value = 111144850
$ws.Range("I93").Value = 166667260
$ws.Range("K93").Value = 166667260
$ws.Range("M93").Value = -166666012
$ws.Range("H100").Value = 202760
$ws.Range("I100").Value = 3450
$ws.Range("K100").Value = 3450
$ws.Range("M100").Value = -2909
$ws.Range("H122").Value = 8287.875
$ws.Range("J122").Value = 7526
$ws.Range("L122").Value = 22578
$ws.Range("N122").Value = -27478
$ws.Range("H132").Value = 11010.593
$ws.Range("I132").Value = 12163.263
$ws.Range("J132").Value = 8273
$ws.Range("K132").Value = 36489.789
$ws.Range("L132").Value = 24819
$ws.Range("M132").Value = -33959.789
$ws.Range("N132").Value = -29879
$ws.Range("H136").Value = 86329.52
$ws.Range("I136").Value = 136432.47
$ws.Range("K136").Value = 409297.41
$ws.Range("M136").Value = -406747.41

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3196.742
$ws.Range("J122").Value = 7840.7144
$ws.Range("L122").Value = 23522.1432
$ws.Range("N122").Value = -28422.1432
$ws.Range("H126").Value = 2508.111
$ws.Range("I126").Value = 2653.7693
$ws.Range("J126").Value = 2129.4
$ws.Range("K126").Value = 7961.3079
$ws.Range("L126").Value = 6388.200000000001
$ws.Range("M126").Value = -5491.3079
$ws.Range("N126").Value = -11328.2
$ws.Range("H132").Value = 43959012
$ws.Range("I132").Value = 55563576
$ws.Range("J132").Value = 2182575
$ws.Range("K132").Value = 166690728
$ws.Range("L132").Value = 6547725
$ws.Range("M132").Value = -166688198
$ws.Range("N132").Value = -6552785
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 4190.933
$ws.Range("I136").Value = 4086.0908
$ws.Range("J136").Value = 4479.25
$ws.Range("K136").Value = 12258.2724
$ws.Range("L136").Value = 13437.75
$ws.Range("M136").Value = -9708.2724
$ws.Range("N136").Value = -18537.75
$ws.Range("H141").Value = 71500
$ws.Range("J141").Value = 71500
$ws.Range("L141").Value = 71500
$ws.Range("N141").Value = -81860
